$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @{
    2 = "Avengers: Age of Ultron "
    3 = "Cinderella "
    4 = "Ant-Man "
    5 = "Do You Believe? "
    6 = "Hot Tub Time Machine 2 "
    7 = "The Water Diviner "
    8 = "Top Five "
    9 = "Shaun the Sheep Movie "
    10 = "Love & Mercy "
    11 = "Far From The Madding Crowd "
    12 = "Black Sea "
    13 = "Leviathan "
    14 = "Unbroken "
    15 = "The Imitation Game "
    16 = "Taken 3 "
    17 = "Ted 2 "
    18 = "Southpaw "
    19 = "Night at the Museum: Secret of the Tomb "
    20 = "Pixels "
    21 = "McFarland, USA "
    22 = "Insidious: Chapter 3 "
    23 = "The Man From U.N.C.L.E. "
    24 = "Run All Night "
    25 = "Trainwreck "
    26 = "Selma "
    27 = "Ex Machina "
    28 = "Still Alice "
    29 = "Wild Tales "
    30 = "The Boy Next Door "
    31 = "Aloha "
    32 = "The Loft "
    33 = "Welcome to Me "
    34 = "Maps to the Stars "
    35 = "Timbuktu "
    36 = "Kingsman: The Secret Service "
    37 = "Tomorrowland "
    38 = "The Divergent Series: Insurgent "
    39 = "Annie "
    40 = "Fantastic Four "
    41 = "Terminator Genisys "
    42 = "Pitch Perfect 2 "
    43 = "Entourage "
    44 = "The Age of Adaline "
    45 = "Hot Pursuit "
    46 = "The DUFF "
    47 = "Project Almanac "
    48 = "Ricki and the Flash "
    49 = "Seventh Son "
    50 = "Mortdecai "
    51 = "Unfinished Business "
    52 = "American Ultra "
    53 = "True Story "
    54 = "Child 44 "
    55 = "Birdman "
    56 = "The Gift "
    57 = "Unfriended "
    58 = "Mr. Turner "
    59 = "American Sniper "
    60 = "Furious 7 "
    61 = "The Hobbit: The Battle of the Five Armies "
    62 = "San Andreas "
    63 = "Straight Outta Compton "
    64 = "Vacation "
    65 = "Chappie "
    66 = "Poltergeist "
    67 = "Paper Towns "
    68 = "Big Eyes "
    69 = "Blackhat "
    70 = "Self/less "
    71 = "Sinister 2 "
    72 = "Little Boy "
    73 = "Me and Earl and The Dying Girl "
    74 = "Maggie "
    75 = "Mad Max: Fury Road "
    76 = "Spy "
    77 = "The SpongeBob Movie: Sponge Out of Water "
    78 = "Paddington "
    79 = "Dope "
    80 = "What We Do in the Shadows "
    81 = "Song of the Sea "
    82 = "Fifty Shades of Grey "
    83 = "Get Hard "
    84 = "Focus "
    85 = "Jupiter Ascending "
    86 = "The Gallows "
    87 = "The Second Best Exotic Marigold Hotel "
    88 = "Strange Magic "
    89 = "The Gunman "
    90 = "Hitman: Agent 47 "
    91 = "Cake "
    92 = "Into the Woods "
    93 = "It Follows "
    94 = "Inherent Vice "
    95 = "A Most Violent Year "
    96 = "While We're Young "
    97 = "Clouds of Sils Maria "
    98 = "Magic Mike XXL "
    99 = "Home "
    100 = "The Wedding Ringer "
    101 = "Woman in Gold "
    102 = "Mission: Impossible â€“ Rogue Nation "
    103 = "Amy "
    104 = "Jurassic World "
    105 = "Minions "
    106 = "Paul Blart: Mall Cop 2 "
    107 = "The Longest Ride "
    108 = "The Lazarus Effect "
    109 = "The Woman In Black 2 Angel of Death "
    110 = "Danny Collins "
    111 = "Inside Out "
    112 = "Mr. Holmes "
    113 = "''71 "
    114 = "Two Days, One Night "
}

foreach ($r in $titles.Keys) {
    $ws.Cells.Item($r, 1).Value = $titles[$r]
}
